$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 corresponds to file 0c38d88b-...
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-08 11:11:57"
$wsZh.Range("G2").Value = "2016-01-08 11:12:41"

# de-de sheet: row 2 corresponds to file 0c38d88b-...
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-08 11:12:07"
$wsDe.Range("G2").Value = "2016-01-08 11:12:57"
